$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CNTF")

# Balance Sheet updates
$ws.Range("D43").Value = 102100   # Net Receivables
$ws.Range("D45").Value = 25800    # Other Current Assets
$ws.Range("D48").Value = 37100    # Property Plant and Equipment
$ws.Range("D52").Value = 149500   # Other Assets
$ws.Range("D58").Value = 8400     # Short/Current Long Term Debt
$ws.Range("D59").Value = 42500    # Other Current Liabilities
$ws.Range("D61").Value = 300      # Long Term Debt

# Cash Flow Statement - Capital Expenditures row
$ws.Range("D91").Value = -40800
$ws.Range("E91").Value = -2100
$ws.Range("F91").Value = -57200
$ws.Range("G91").Value = -49900
$ws.Range("H91").Value = -9900
$ws.Range("I91").Value = -15700
